# Końcowe wyniki wyzarzanie w tabeli
# Fill in the missing "100" row (row 5) results on the "127" sheet with the
# final annealing measurements, and move the active selection from C22 to C19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("127")
$ws.Activate()

# Row 5 (parameter value = 100) was still empty for columns C, D and E.
# Column E uses a "Text" (@) number format, so we briefly switch it to a
# numeric format while writing the value, then restore it, to make sure the
# value is stored as a real number (matching the other rows) instead of text.
$ws.Cells.Item(5, 3).Value = 118532.858643336
$ws.Cells.Item(5, 4).Value = 121139.51

$eFmt = $ws.Cells.Item(5, 5).NumberFormat
$ws.Cells.Item(5, 5).NumberFormat = "0.00"
$ws.Cells.Item(5, 5).Value = 6856.91
$ws.Cells.Item(5, 5).NumberFormat = $eFmt

# Move the selection, as left by the author after finishing the table.
$ws.Range("C19").Select()
